$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "3C" -> "3c" : rename class label everywhere it appears (rows that keep a 3-class
# value in column C). Row 44 is fully rewritten below with a different class, so it
# is excluded here.
$classCells = @("C2","C6","C9","C14","C20","C22","C24","C28","C32","C33","C38","C39","C41")
foreach ($cellRef in $classCells) {
    $ws.Range($cellRef).Value = "3c"
}

# Row 44: was Stephan/Fuchs/3C/18.12.2023 08:00/18.12.2023 09:40/100
#  -> now Max/Schmitz/4a/18.12.2023 08:00/18.12.2023 16:19/499  (Abmeldung moved to afternoon check-in window)
$ws.Cells.Item(44, 1).Value = "Max"
$ws.Cells.Item(44, 2).Value = "Schmitz"
$ws.Cells.Item(44, 3).Value = "4a"
$ws.Cells.Item(44, 4).Value = "18.12.2023 08:00"
$ws.Cells.Item(44, 5).Value = "18.12.2023 16:19"
$ws.Cells.Item(44, 6).Value = 499

# Row 45: was Max/Schmitz/4a/18.12.2023 08:00/18.12.2023 08:01/1
#  -> now Eli/Enders/4a/18.12.2023 09:53/18.12.2023 16:20/387
$ws.Cells.Item(45, 1).Value = "Eli"
$ws.Cells.Item(45, 2).Value = "Enders"
$ws.Cells.Item(45, 3).Value = "4a"
$ws.Cells.Item(45, 4).Value = "18.12.2023 09:53"
$ws.Cells.Item(45, 5).Value = "18.12.2023 16:20"
$ws.Cells.Item(45, 6).Value = 387

# Row 46: was Detlef/Soost/1a/18.12.2023 08:00/18.12.2023 08:01/1
#  -> now Stephan/Fuchs/3c/20.12.2023 15:53/20.12.2023 15:53/0 (Anmeldung only, no Abmeldung yet -> 0 duration)
$ws.Cells.Item(46, 1).Value = "Stephan"
$ws.Cells.Item(46, 2).Value = "Fuchs"
$ws.Cells.Item(46, 3).Value = "3c"
$ws.Cells.Item(46, 4).Value = "20.12.2023 15:53"
$ws.Cells.Item(46, 5).Value = "20.12.2023 15:53"
$ws.Cells.Item(46, 6).Value = 0

# Row 47 (new): Stephan/Fuchs/3c/15.02.2024 13:43/15.02.2024 14:04/21
$ws.Cells.Item(47, 1).Value = "Stephan"
$ws.Cells.Item(47, 2).Value = "Fuchs"
$ws.Cells.Item(47, 3).Value = "3c"
$ws.Cells.Item(47, 4).Value = "15.02.2024 13:43"
$ws.Cells.Item(47, 5).Value = "15.02.2024 14:04"
$ws.Cells.Item(47, 6).Value = 21
